# Amélioration sur conditions de stockage et date de naissance
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Fill in patient name data (columns A=Nom, B=nom de jeune fille, C=Prénom) ---
# Filled column by column (A first, then B, then C) so the shared-string table
# is built up in the same order as the source data.
$ws.Range("A2").Value = "MARTIN"
$ws.Range("A3").Value = "BERNARD"

$ws.Range("B2").Value = "ANDRE"

$ws.Range("C2").Value = "Marie"
$ws.Range("C3").Value = "Jean"

# --- 2. Normalize formatting of the Nom / nom de jeune fille / Prénom columns ---
# Columns A:C for every data row (2-41) get the same look as the rest of the
# row data (centered, same date-like numeric formatting as column D).
$ws.Range("D2").Copy()
$ws.Range("A2:C41").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- 3. Update "sexe" (U column) values for the data rows ---
$ws.Range("U4").Value = "M"
$ws.Range("U6").Value = "M"
$ws.Range("U10").Value = "M"
$ws.Range("U11").Value = "F"
$ws.Range("U12").Value = "F"
$ws.Range("U13").Value = "F"
$ws.Range("U14").Value = "F"
$ws.Range("U15").Value = "M"
$ws.Range("U16").Value = "M"
$ws.Range("U17").Value = "F"
$ws.Range("U18").Value = "M"
$ws.Range("U19").Value = "M"
$ws.Range("U22").Value = "M"
$ws.Range("U23").Value = "M"
$ws.Range("U24").Value = "F"
$ws.Range("U26").Value = "F"
$ws.Range("U27").Value = "M"
$ws.Range("U30").Value = "M"
$ws.Range("U33").Value = "M"
$ws.Range("U39").Value = "M"
$ws.Range("U40").Value = "M"

# Rows 29, 31 and 38 additionally carried a stray red-font style on the
# "sexe" cell - bring them in line with the rest of the column first, then
# set their values.
$ws.Range("U2").Copy()
$ws.Range("U29").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("U31").PasteSpecial(-4122)
$ws.Range("U38").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("U29").Value = "M"
$ws.Range("U31").Value = "M"
$ws.Range("U38").Value = "M"

# --- 4. Remove the trailing blank rows (42-46) ---
$ws.Rows("42:46").Delete()

# --- 5. Restore the selection to match the saved view ---
[void]$ws.Range("G7").Select()
